# Scaling updates to better match EMEP, etc.
# - "year" sheet: add a new "select_scaling_year" column, tweak the scaling
#   window for the existing ("mkd") row, and add a new "fin" row.
# - "map" sheet: carry over the author's new cursor/selection position.

$wb = $excel.ActiveWorkbook

$wsMap    = $wb.Worksheets.Item("map")
$wsMethod = $wb.Worksheets.Item("method")
$wsYear   = $wb.Worksheets.Item("year")

# ---------------------------------------------------------------------------
# "year" sheet data edits
# ---------------------------------------------------------------------------

# Insert a new column H ("select_scaling_year"); the former column H
# ("Comment") slides over to column I.
$oldH1 = $wsYear.Cells.Item(1, 8).Value2
$wsYear.Cells.Item(1, 9).Value = $oldH1
$wsYear.Cells.Item(1, 8).Value = "select_scaling_year"

# Row 2 ("mkd"/"all"): scale from 2000 instead of 1990, clear the old
# comment (now "NA"), and record the new rationale in column I.
$wsYear.Cells.Item(2, 6).Value = 2000
$wsYear.Cells.Item(2, 8).Value = "NA"
$wsYear.Cells.Item(2, 9).Value = "Scale from 2000 so as to be closer to EMEP trend"

# New row 3: "fin"/"all" scaling window 1982-2020, with its own comment.
$wsYear.Cells.Item(3, 1).Value = "fin"
$wsYear.Cells.Item(3, 2).Value = "all"
$wsYear.Cells.Item(3, 3).Value = "NA"
$wsYear.Cells.Item(3, 4).Value = "NA"
$wsYear.Cells.Item(3, 5).Value = "NA"
$wsYear.Cells.Item(3, 6).Value = 1982
$wsYear.Cells.Item(3, 7).Value = 2020
$wsYear.Cells.Item(3, 8).Value = "NA"
$wsYear.Cells.Item(3, 9).Value = "Don't scale 1981 to avoid reporting mistake in inventory"

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------

# "map" sheet cursor moved to B37.
$wsMap.Range("B37").Select()

# "year" sheet selection now covers the full new row 3.
$wsYear.Select()
$wsYear.Range("A3:XFD3").Select()
